$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.001754667048134761
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 10137753.70137369
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 10209273.24794223
